$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F ("想去人数") for the listed rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 177
$ws1.Range("F5").Value = 192
$ws1.Range("F7").Value = 1135
$ws1.Range("F8").Value = 379
$ws1.Range("F9").Value = 191
$ws1.Range("F12").Value = 373
$ws1.Range("F13").Value = 392
$ws1.Range("F14").Value = 784
$ws1.Range("F15").Value = 172
$ws1.Range("F16").Value = 720
$ws1.Range("F17").Value = 282
$ws1.Range("F18").Value = 78
$ws1.Range("F19").Value = 1003
$ws1.Range("F20").Value = 455
$ws1.Range("F21").Value = 261
$ws1.Range("F22").Value = 82
$ws1.Range("F23").Value = 379
$ws1.Range("F25").Value = 40
$ws1.Range("F26").Value = 468

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 364
$ws2.Range("F7").Value = 282

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 344

# Sheet "全部类型" (sheet4) - combined listing of all events
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 344
$ws4.Range("F5").Value = 177
$ws4.Range("F7").Value = 192
$ws4.Range("F9").Value = 1135
$ws4.Range("F10").Value = 379
$ws4.Range("F11").Value = 191
$ws4.Range("F14").Value = 364
$ws4.Range("F17").Value = 373
$ws4.Range("F19").Value = 282
$ws4.Range("F20").Value = 392
$ws4.Range("F21").Value = 784
$ws4.Range("F22").Value = 172
$ws4.Range("F23").Value = 720
$ws4.Range("F24").Value = 282
$ws4.Range("F25").Value = 78
$ws4.Range("F26").Value = 1003
$ws4.Range("F27").Value = 455
$ws4.Range("F30").Value = 261
$ws4.Range("F31").Value = 82
$ws4.Range("F32").Value = 379
$ws4.Range("F36").Value = 40
$ws4.Range("F38").Value = 468
